$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11

$cell = $t.Cell(1, 1)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "72 x 96" + $vt + "  9    6" + $vt + "  ----" + $vt + "7|    |" + $vt + "2|    |"

$cell = $t.Cell(1, 2)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "73 x 75" + $vt + "  7    5" + $vt + "  ----" + $vt + "7|    |" + $vt + "3|    |"

$cell = $t.Cell(1, 3)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "95 x 37" + $vt + "  3    7" + $vt + "  ----" + $vt + "9|    |" + $vt + "5|    |"

$cell = $t.Cell(2, 1)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "24 x 41" + $vt + "  4    1" + $vt + "  ----" + $vt + "2|    |" + $vt + "4|    |"

$cell = $t.Cell(2, 2)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "51 x 97" + $vt + "  9    7" + $vt + "  ----" + $vt + "5|    |" + $vt + "1|    |"

$cell = $t.Cell(2, 3)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "17 x 49" + $vt + "  4    9" + $vt + "  ----" + $vt + "1|    |" + $vt + "7|    |"

$cell = $t.Cell(3, 1)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "97 x 33" + $vt + "  3    3" + $vt + "  ----" + $vt + "9|    |" + $vt + "7|    |"

$cell = $t.Cell(3, 2)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "87 x 47" + $vt + "  4    7" + $vt + "  ----" + $vt + "8|    |" + $vt + "7|    |"

$cell = $t.Cell(3, 3)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "86 x 86" + $vt + "  8    6" + $vt + "  ----" + $vt + "8|    |" + $vt + "6|    |"

$cell = $t.Cell(4, 1)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "62 x 77" + $vt + "  7    7" + $vt + "  ----" + $vt + "6|    |" + $vt + "2|    |"

$cell = $t.Cell(4, 2)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "62 x 85" + $vt + "  8    5" + $vt + "  ----" + $vt + "6|    |" + $vt + "2|    |"

$cell = $t.Cell(4, 3)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "45 x 19" + $vt + "  1    9" + $vt + "  ----" + $vt + "4|    |" + $vt + "5|    |"

$cell = $t.Cell(5, 1)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "18 x 24" + $vt + "  2    4" + $vt + "  ----" + $vt + "1|    |" + $vt + "8|    |"

$cell = $t.Cell(5, 2)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "83 x 21" + $vt + "  2    1" + $vt + "  ----" + $vt + "8|    |" + $vt + "3|    |"

$cell = $t.Cell(5, 3)
$r2 = $cell.Range
$r2.End = $r2.End - 1
$r2.Text = "73 x 96" + $vt + "  9    6" + $vt + "  ----" + $vt + "7|    |" + $vt + "3|    |"

Write-Host "Done updating lattice multiplication cells"

